# Append newly-tracked national park names ("G" and "H" entries) to the
# Park Codes worksheet. The sheet already lists parks alphabetically in
# column A (shared strings), one per row, with a blank separator row
# between each starting letter group (e.g. rows 33, 63, 123, 135, 150 in
# the original data separate A/B, B/C, C/D, D/E, E/F). This change adds
# the "G" park names directly below the existing data (row 191 stays
# blank, continuing that separator pattern), a blank separator row at
# 226 for G/H, then the "H" park names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "G" parks (rows 192-225); row 226 is left blank as the G/H separator.
$ws.Range("A192").Value = 'Gates Of The Arctic National Park & Preserve'
$ws.Range("A193").Value = 'Gateway Arch National Park'
$ws.Range("A194").Value = 'Gateway National Recreation Area'
$ws.Range("A195").Value = 'Gauley River National Recreation Area'
$ws.Range("A196").Value = 'General Grant National Memorial'
$ws.Range("A197").Value = 'George Rogers Clark National Historical Park'
$ws.Range("A198").Value = 'George Washington Birthplace National Monument'
$ws.Range("A199").Value = 'George Washington Carver National Monument'
$ws.Range("A200").Value = 'George Washington Memorial Parkway'
$ws.Range("A201").Value = 'Gettysburg National Military Park'
$ws.Range("A202").Value = 'Gila Cliff Dwellings National Monument'
$ws.Range("A203").Value = 'Glacier Bay National Park & Preserve'
$ws.Range("A204").Value = 'Glacier National Park'
$ws.Range("A205").Value = 'Glen Canyon National Recreation Area'
$ws.Range("A206").Value = 'Glen Echo Park'
$ws.Range("A207").Value = 'Gloria Dei Church National Historic Site'
$ws.Range("A208").Value = 'Golden Gate National Recreation Area'
$ws.Range("A209").Value = 'Golden Spike National Historical Park'
$ws.Range("A210").Value = 'Governors Island National Monument'
$ws.Range("A211").Value = 'Grand Canyon National Park'
$ws.Range("A212").Value = 'Grand Portage National Monument'
$ws.Range("A213").Value = 'Grand Teton National Park'
$ws.Range("A214").Value = 'Grant-Kohrs Ranch National Historic Site'
$ws.Range("A215").Value = 'Great Basin National Park'
$ws.Range("A216").Value = 'Great Egg Harbor River'
$ws.Range("A217").Value = 'Great Falls Park'
$ws.Range("A218").Value = 'Great Sand Dunes National Park & Preserve'
$ws.Range("A219").Value = 'Great Smoky Mountains National Park'
$ws.Range("A220").Value = 'Green Springs'
$ws.Range("A221").Value = 'Greenbelt Park'
$ws.Range("A222").Value = 'Guadalupe Mountains National Park'
$ws.Range("A223").Value = 'Guilford Courthouse National Military Park'
$ws.Range("A224").Value = 'Gulf Islands National Seashore'
$ws.Range("A225").Value = 'Gullah/Geechee Cultural Heritage Corridor'
$ws.Range("A227").Value = 'Hagerman Fossil Beds National Monument'
$ws.Range("A228").Value = 'Haleakalā National Park'
$ws.Range("A229").Value = 'Hamilton Grange National Memorial'
$ws.Range("A230").Value = 'Hampton National Historic Site'
$ws.Range("A231").Value = 'Harmony Hall'
$ws.Range("A232").Value = 'Harpers Ferry National Historical Park'
$ws.Range("A233").Value = 'Harriet Tubman National Historical Park'
$ws.Range("A234").Value = 'Harriet Tubman Underground Railroad National Historical Park'
$ws.Range("A235").Value = 'Harry S Truman National Historic Site'
$ws.Range("A236").Value = 'Hawai''i Volcanoes National Park'
$ws.Range("A237").Value = 'Herbert Hoover National Historic Site'
$ws.Range("A238").Value = 'Historic Jamestowne Part of Colonial National Historical Park'
$ws.Range("A239").Value = 'Home Of Franklin D Roosevelt National Historic Site'
$ws.Range("A240").Value = 'Homestead National Monument of America'
$ws.Range("A241").Value = 'Honouliuli National Historic Site'
$ws.Range("A242").Value = 'Hopewell Culture National Historical Park'
$ws.Range("A243").Value = 'Hopewell Furnace National Historic Site'
$ws.Range("A244").Value = 'Horseshoe Bend National Military Park'
$ws.Range("A245").Value = 'Hot Springs National Park'
$ws.Range("A246").Value = 'Hovenweep National Monument'
$ws.Range("A247").Value = 'Hubbell Trading Post National Historic Site'
$ws.Range("A248").Value = 'Hudson River Valley National Heritage Area'

# Scroll the view down to show the newly added rows and select the last
# populated cell, matching where the author's cursor ended up.
$excel.ActiveWindow.ScrollRow = 233
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A248").Select() | Out-Null
